$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the four query strings: "name" column -> "user" column (both the
# selected column and the WHERE clause column reference).
$ws.Range("M2").Value = "select host,user,AUTHENTICATION_STRING from mysql.user where user='MYDCLTEST1'"
$ws.Range("M3").Value = "select host,user,AUTHENTICATION_STRING from mysql.user where user='mydcltest2'"
$ws.Range("M4").Value = "select host,user,AUTHENTICATION_STRING from mysql.user where user='MydclTest3'"
$ws.Range("M5").Value = "select user,AUTHENTICATION_STRING from mysql.user where user='mydcltest4'"

# Update the sheet view: drop the topLeftCell scroll position and move the
# active selection to M17.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M17").Select()
